# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Values are written with a leading apostrophe to force text entry (matching
# the sheet's existing inline-string cells), then the style is reset back to
# "Normal" so the quote-prefix formatting doesn't linger on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''41.774.33'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +0.55%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = '''2.231.20'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -0.96%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = '''  +0.00%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = '''231.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -1.24%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = '''0.621'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -2.44%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = '''59.99'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  -7.47%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = '''  -0.03%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = '''0.402'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -1.86%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = '''57.86'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -2.78%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = '''0.0895'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -0.76%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("E12").Value = '''  -0.99%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = '''2.562.18'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -1.06%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = '''15.41'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -4.82%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = '''22.50'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +0.24%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = '''  -0.57%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = '''0.799'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -4.31%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = '''2.241.34'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -0.58%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = '''41.693.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +0.60%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = '''0.0₃0909'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -1.07%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = '''72.41'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -2.31%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = '''6.12'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -1.23%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = '''247.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -2.21%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = '''0.998'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -0.18%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = '''  -2.13%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = '''  -0.61%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = '''9.75'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -0.85%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = '''169.41'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -2.20%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = '''0.140'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  -2.73%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = '''19.89'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -3.01%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = '''  -2.52%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = '''  -8.29%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = '''  -2.18%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = '''4.98'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -0.08%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = '''4.68'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -1.36%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = '''  +3.25%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = '''6.54'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -9.24%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = '''2.40'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -2.65%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = '''3.60'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -7.74%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = '''0.000242'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +2.78%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D42").Value = '''0.0239'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +0.57%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = '''8.64'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -1.97%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = '''  -1.54%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = '''98.84'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -3.40%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = '''0.0964'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +2.18%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = '''Maker'
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = '''https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = '''1.472.66'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -2.85%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = '''FTXToken'
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = '''https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = '''4.35'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -10.69%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = '''16.57'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -8.90%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = '''  -1.74%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = '''  -3.07%  '
$ws.Range("E51").Style = "Normal"
